$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.845.31"
$ws.Range("E2").Value = "  -4.14%  "
$ws.Range("D3").Value = "3.658.51"
$ws.Range("E3").Value = "  -3.83%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'590.38"
$ws.Range("E5").Value = "  -3.92%  "
$ws.Range("D6").Value = "'163.36"
$ws.Range("E6").Value = "  -7.83%  "
$ws.Range("D7").Value = "3.657.10"
$ws.Range("E7").Value = "  -3.84%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "'0.519"
$ws.Range("E9").Value = "  -1.97%  "
$ws.Range("D10").Value = "'0.158"
$ws.Range("E10").Value = "  -5.92%  "
$ws.Range("D11").Value = "'6.08"
$ws.Range("E11").Value = "  -6.27%  "
$ws.Range("D12").Value = "'0.457"
$ws.Range("E12").Value = "  -5.71%  "
$ws.Range("D13").Value = "'37.12"
$ws.Range("E13").Value = "  -6.77%  "
$ws.Range("D14").Value = "'0.0000237"
$ws.Range("E14").Value = "  -6.91%  "
$ws.Range("D15").Value = "4.270.85"
$ws.Range("E15").Value = "  -3.78%  "
$ws.Range("D16").Value = "3.662.33"
$ws.Range("E16").Value = "  -3.77%  "
$ws.Range("D17").Value = "66.980.12"
$ws.Range("E17").Value = "  -4.04%  "
$ws.Range("E18").Value = "  -4.43%  "
$ws.Range("D19").Value = "'7.06"
$ws.Range("E19").Value = "  -6.60%  "
$ws.Range("D20").Value = "'16.83"
$ws.Range("E20").Value = "  +0.79%  "
$ws.Range("D21").Value = "'485.56"
$ws.Range("E21").Value = "  -4.59%  "
$ws.Range("D22").Value = "'8.98"
$ws.Range("E22").Value = "  -6.44%  "
$ws.Range("D23").Value = "'0.707"
$ws.Range("E23").Value = "  -4.06%  "
$ws.Range("D24").Value = "'84.82"
$ws.Range("E24").Value = "  -1.86%  "
$ws.Range("D25").Value = "'2.26"
$ws.Range("E25").Value = "  -8.47%  "
$ws.Range("D26").Value = "'0.0000137"
$ws.Range("E26").Value = "  -5.73%  "
$ws.Range("D27").Value = "'11.98"
$ws.Range("E27").Value = "  -5.91%  "
$ws.Range("D28").Value = "'0.994"
$ws.Range("E28").Value = "  -0.57%  "
$ws.Range("D29").Value = "'9.82"
$ws.Range("E29").Value = "  -6.91%  "
$ws.Range("D30").Value = "'2.89"
$ws.Range("E30").Value = "  -2.94%  "
$ws.Range("D31").Value = "'2.33"
$ws.Range("E31").Value = "  -7.40%  "
$ws.Range("D32").Value = "'7.62"
$ws.Range("E32").Value = "  -5.45%  "
$ws.Range("D33").Value = "'31.27"
$ws.Range("E33").Value = "  -1.08%  "
$ws.Range("D34").Value = "3.797.96"
$ws.Range("E34").Value = "  -3.84%  "
$ws.Range("D35").Value = "3.596.59"
$ws.Range("E35").Value = "  -3.79%  "
$ws.Range("D36").Value = "'0.105"
$ws.Range("E36").Value = "  -7.69%  "
$ws.Range("D37").Value = "'1.00"
$ws.Range("E37").Value = "  +0.05%  "
$ws.Range("D38").Value = "'0.987"
$ws.Range("E38").Value = "  -5.55%  "
$ws.Range("D39").Value = "'5.68"
$ws.Range("E39").Value = "  -7.22%  "
$ws.Range("D40").Value = "'0.130"
$ws.Range("E40").Value = "  -7.81%  "
$ws.Range("D41").Value = "'0.318"
$ws.Range("E41").Value = "  -6.07%  "
$ws.Range("D42").Value = "'432.47"
$ws.Range("E42").Value = "  -10.04%  "
$ws.Range("D43").Value = "'48.43"
$ws.Range("E43").Value = "  -2.69%  "
$ws.Range("D44").Value = "'1.89"
$ws.Range("E44").Value = "  -8.14%  "
$ws.Range("D45").Value = "'2.74"
$ws.Range("E45").Value = "  -9.09%  "
$ws.Range("D46").Value = "'8.25"
$ws.Range("E46").Value = "  -3.93%  "
$ws.Range("E47").Value = "  -0.03%  "
$ws.Range("D48").Value = "'141.79"
$ws.Range("E48").Value = "  +1.93%  "
$ws.Range("D49").Value = "'39.39"
$ws.Range("E49").Value = "  -10.85%  "
$ws.Range("D50").Value = "2.737.50"
$ws.Range("E50").Value = "  -6.79%  "
$ws.Range("D51").Value = "'0.0343"
$ws.Range("E51").Value = "  -5.56%  "
